# Project "Sample Project" save: update Rule R40's row (row 11) so the
# "Rule" column (B11) now reads "1" instead of "R40".
#
# The target cell must hold the literal text "1" (not the number 1), so
# the cell is pre-formatted as Text ("@") before the value is written;
# otherwise Excel would auto-coerce the numeric-looking literal into a
# number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
